$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the STATUS column (G) that need to be marked "DONE" -- this
# includes every task/subtask row (the blank separator rows 5,10,14,19,
# 25,29,33 are intentionally left untouched).
$rows = @(2,3,4, 6,7,8,9, 11,12,13, 15,16,17,18, 20,21,22,23,24, 26,27,28, 30,31,32)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $cell.Value = "DONE"
    $cell.Font.Name = "Bookman Old Style"
    $cell.Font.Size = 14
    $cell.Font.Bold = $true
    $cell.Font.Color = 255
    $cell.HorizontalAlignment = -4108
}

# Update the sheet view: zoom level and active selection.
$excel.ActiveWindow.Zoom = 39
$ws.Range("G23:G24").Select() | Out-Null
